# Changes of 6th May 2022
#
# Update the tracking-number values in column P ("ActualRate"/shipment
# tracking column) on rows 2 and 4 of the CRUD-operation test sheet to
# the new shipment tracking numbers. The values must remain stored as
# text (shared-string) cells, exactly like the pre-existing tracking
# numbers in that column, with no incidental style/number-format change
# on the cell.
#
# Simply assigning a numeric-looking string to .Value lets Excel infer a
# numeric type (losing the text/shared-string representation), and
# forcing text via NumberFormat="@" leaves a permanent (if unused) extra
# cell-style entry behind. Instead, compute the text via TEXT() and
# paste back only the resulting value, which preserves the original
# "no explicit style" cell shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("P2").Formula = "=TEXT(320018475104,""0"")"
$ws.Range("P2").Copy()
$ws.Range("P2").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("P4").Formula = "=TEXT(320018475115,""0"")"
$ws.Range("P4").Copy()
$ws.Range("P4").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = 0
